$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the GitHub Usernames column (C) for the existing teams, plus the
# new team in row 5, as one block (matches original authoring order)
$ws.Range("C2").Value = "inba-m06, KeerthanaGb14, Lokesh-1511"
$ws.Range("C3").Value = "Gokul-K-19, Sriram-2705, arumugaperumal06"
$ws.Range("C4").Value = "Jayasuriya-ENGINEER, Thulasi-Ravichandran, vishy-18"
$ws.Range("C5").Value = "Abishree730, SubhasiniPaulpandi23, Praveen95-cs"

# New TeamID column entries
$ws.Range("A5").Value = "T004"
$ws.Range("A6").Value = "T005"

# New Team Name column entries
$ws.Range("B5").Value = "FGH"
$ws.Range("B6").Value = "HJI"

# GitHub usernames for the second new team
$ws.Range("C6").Value = "rakeshmm05, VethavalliGM"

# PS ID / PS for row 5
$ws.Range("D5").Value = "L04"
$ws.Range("E5").Value = "CLOUD "

# PS ID / PS for row 6
$ws.Range("D6").Value = "L05"
$ws.Range("E6").Value = "NETWORKS"

# Update selection to mirror the final cursor position
$ws.Range("E6").Select()
